$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Remove the stray "sdsd" / "sdsds" cells in column H (rows 2 and 3) ---
# Use Clear() (not ClearContents()) so the now-empty cells drop out of the
# saved XML entirely, matching the diff which removes the <c> elements.
$ws.Range("H2").Clear() | Out-Null
$ws.Range("H3").Clear() | Out-Null

# --- Row 4: "assurance" word-form / synonym / IPA / meaning / definition / sentence ---
# Fill in the same order the shared-strings table in the target file implies
# (new-word text, then IPA, then synonym, then meaning, then definition last)
$ws.Range("B4").Value2 = "n"
$ws.Range("G4").Value2 = "The sales assurance that the missing keyboard would be replaced the next day."
$ws.Range("D4").Value2 = "/əˈʃʊərəns/"
$ws.Range("C4").Value2 = "aguarantee; confidence, promise"
$ws.Range("E4").Value2 = "Đảm bảo, tự tin"

# --- Row 5: "cancellation" word-form / synonym / IPA / meaning / definition / sentence ---
$ws.Range("B5").Value2 = "n"
$ws.Range("D5").Value2 = "/ˌkænsəˈleɪʃn/"
$ws.Range("E5").Value2 = "Sự hủy bỏ"
$ws.Range("C5").Value2 = "annulment, stopping"
$ws.Range("G5").Value2 = "The cancellation of hẻ flight caused her problems for the rest of the week."

# --- Definitions added last for both new rows ---
$ws.Range("F4").Value2 = "a statement that something will certainly be true or will certainly happen, particularly when there has been doubt about it"
$ws.Range("F5").Value2 = "a decision to stop something that has already been arranged from happening; a statement that something will not happen"

# --- Column width adjustments (columns C and F got wider) ---
$ws.Columns.Item(3).ColumnWidth = 32.42
$ws.Columns.Item(6).ColumnWidth = 39.75

# --- Row 2 is now shorter since its tall wrapped cell no longer forces extra height ---
$ws.Rows.Item(2).RowHeight = 33

# --- Selection moved to B6 ---
$ws.Range("B6").Select() | Out-Null
